$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header columns (AD, AE, AF) - "Wins", "Losses", "Ties"
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting of the existing header cells (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every player data row
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 88
    $ws.Cells.Item($r, 31).Value = 74
    $ws.Cells.Item($r, 32).Value = 0
}
